$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update October 2025 row (row 23) stats
$ws.Range("B23").Value = 6328
$ws.Range("D23").Value = 5892647
$ws.Range("E23").Value = 931.2021175726928
$ws.Range("F23").Value = 8.5792724776939
$ws.Range("H23").Value = 26.26434218936586
